$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that can look numeric (e.g. "209.43").
# Excel would normally coerce these into floating point numbers on
# assignment, which both loses the exact text representation (e.g.
# trailing zeros / precision) and changes the stored cell type from
# text to number. To preserve the original text semantics we briefly
# force a text number format before writing the value, then restore
# the cell to the default "Normal" style so no visible formatting
# change is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.926.75'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.98%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.631.16'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.62%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5201'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.00%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -3.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06220'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.20'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.29%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07558'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.635.90'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.344'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.859.48'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5414'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅7920'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.17%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.50'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.44%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.911.99'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.25%  '
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.606'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '184.33'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.79%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.01'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.056'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.96%  '
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.64'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.75%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1208'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.98%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.329'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.94%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.48'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.367'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05909'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.99%  '
$ws.Range("E31").Value = '  -3.46%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.344'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.348'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.601'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.76%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9687'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.384'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.738'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5766'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.87%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01594'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8343'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.606'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.91%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.010.82'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.63%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.57'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.38%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.783.95'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₈109'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.91%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9967'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.18'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.26%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.941'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05171'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4223'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.79%  '
